$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 192; everything below (old 192:233) shifts
# down to 194:235, keeping the weekly series intact with a new week's data
# (date 44588) inserted at the top of this block.
$ws.Rows("192:193").Insert()

# Row 192: Betarraga "Primera" quality for the new week (2022-01-27 / 44588)
$ws.Cells.Item(192, 1).Value = 8
$ws.Cells.Item(192, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(192, 3).Value = "Coquimbo"
$ws.Cells.Item(192, 4).Value = 44588
$ws.Cells.Item(192, 5).Value = 4
$ws.Cells.Item(192, 6).Value = 100114014
$ws.Cells.Item(192, 7).Value = "Betarraga"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 2800
$ws.Cells.Item(192, 11).Value = 450
$ws.Cells.Item(192, 12).Value = 500
$ws.Cells.Item(192, 13).Value = 475
$ws.Cells.Item(192, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(192, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(192, 16).Value = 158
$ws.Cells.Item(192, 17).Value = 3
$ws.Cells.Item(192, 18).Value = "Hortaliza"

# Row 193: Betarraga "Segunda" quality for the same new week
$ws.Cells.Item(193, 1).Value = 8
$ws.Cells.Item(193, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(193, 3).Value = "Coquimbo"
$ws.Cells.Item(193, 4).Value = 44588
$ws.Cells.Item(193, 5).Value = 4
$ws.Cells.Item(193, 6).Value = 100114014
$ws.Cells.Item(193, 7).Value = "Betarraga"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Segunda"
$ws.Cells.Item(193, 10).Value = 1600
$ws.Cells.Item(193, 11).Value = 350
$ws.Cells.Item(193, 12).Value = 400
$ws.Cells.Item(193, 13).Value = 375
$ws.Cells.Item(193, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(193, 15).Value = "Provincia del Elqu" + [char]0x00ED
$ws.Cells.Item(193, 16).Value = 125
$ws.Cells.Item(193, 17).Value = 3
$ws.Cells.Item(193, 18).Value = "Hortaliza"
